$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$updates = @(
    @{ Cell = 'D2'; Value = '42.087.17' }
    @{ Cell = 'E2'; Value = '  +0.67%  ' }
    @{ Cell = 'D3'; Value = '2.264.19' }
    @{ Cell = 'E3'; Value = '  -0.25%  ' }
    @{ Cell = 'D4'; Value = '1.00' }
    @{ Cell = 'E4'; Value = '  +0.02%  ' }
    @{ Cell = 'D5'; Value = '305.60' }
    @{ Cell = 'E5'; Value = '  +0.16%  ' }
    @{ Cell = 'D6'; Value = '95.77' }
    @{ Cell = 'E6'; Value = '  +3.04%  ' }
    @{ Cell = 'E8'; Value = '  -0.02%  ' }
    @{ Cell = 'D9'; Value = '0.490' }
    @{ Cell = 'E9'; Value = '  +0.92%  ' }
    @{ Cell = 'D10'; Value = '35.15' }
    @{ Cell = 'E10'; Value = '  +7.62%  ' }
    @{ Cell = 'E11'; Value = '  -1.09%  ' }
    @{ Cell = 'E12'; Value = '  -0.20%  ' }
    @{ Cell = 'E13'; Value = '  -0.47%  ' }
    @{ Cell = 'D14'; Value = '2.616.00' }
    @{ Cell = 'E14'; Value = '  -0.23%  ' }
    @{ Cell = 'D16'; Value = '2.271.36' }
    @{ Cell = 'E16'; Value = '  -0.02%  ' }
    @{ Cell = 'D17'; Value = '0.794' }
    @{ Cell = 'E17'; Value = '  +0.98%  ' }
    @{ Cell = 'D18'; Value = '41.981.19' }
    @{ Cell = 'E18'; Value = '  +0.53%  ' }
    @{ Cell = 'E19'; Value = '  -4.20%  ' }
    @{ Cell = 'D20'; Value = '0.0₃0904' }
    @{ Cell = 'E20'; Value = '  -0.53%  ' }
    @{ Cell = 'E21'; Value = '  -0.01%  ' }
    @{ Cell = 'D22'; Value = '67.64' }
    @{ Cell = 'E22'; Value = '  -0.44%  ' }
    @{ Cell = 'D23'; Value = '238.00' }
    @{ Cell = 'E23'; Value = '  -2.33%  ' }
    @{ Cell = 'E24'; Value = '  -0.64%  ' }
    @{ Cell = 'B25'; Value = 'ImmutableX' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D25'; Value = '1.93' }
    @{ Cell = 'E25'; Value = '  -0.11%  ' }
    @{ Cell = 'B26'; Value = 'Dai' }
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D26'; Value = '1.00' }
    @{ Cell = 'E26'; Value = '  -0.05%  ' }
    @{ Cell = 'D27'; Value = '23.71' }
    @{ Cell = 'E27'; Value = '  -1.22%  ' }
    @{ Cell = 'D28'; Value = '36.63' }
    @{ Cell = 'E28'; Value = '  +5.48%  ' }
    @{ Cell = 'D29'; Value = '9.51' }
    @{ Cell = 'E29'; Value = '  -1.21%  ' }
    @{ Cell = 'E30'; Value = '  +1.85%  ' }
    @{ Cell = 'D31'; Value = '159.96' }
    @{ Cell = 'E31'; Value = '  +0.24%  ' }
    @{ Cell = 'E32'; Value = '  -1.93%  ' }
    @{ Cell = 'E33'; Value = '  +0.01%  ' }
    @{ Cell = 'D34'; Value = '3.19' }
    @{ Cell = 'E34'; Value = '  +5.32%  ' }
    @{ Cell = 'D35'; Value = '0.0739' }
    @{ Cell = 'E35'; Value = '  -0.45%  ' }
    @{ Cell = 'D36'; Value = '17.10' }
    @{ Cell = 'E36'; Value = '  +0.39%  ' }
    @{ Cell = 'E37'; Value = '  +0.50%  ' }
    @{ Cell = 'E38'; Value = '  -0.81%  ' }
    @{ Cell = 'E39'; Value = '  +2.09%  ' }
    @{ Cell = 'E40'; Value = '  -1.72%  ' }
    @{ Cell = 'E41'; Value = '  +2.19%  ' }
    @{ Cell = 'E42'; Value = '  +7.14%  ' }
    @{ Cell = 'D43'; Value = '1.982.17' }
    @{ Cell = 'E43'; Value = '  -1.50%  ' }
    @{ Cell = 'D44'; Value = '19.05' }
    @{ Cell = 'E44'; Value = '  -3.43%  ' }
    @{ Cell = 'E45'; Value = '  +0.29%  ' }
    @{ Cell = 'D46'; Value = '2.93' }
    @{ Cell = 'E46'; Value = '  +0.32%  ' }
    @{ Cell = 'D47'; Value = '9.93' }
    @{ Cell = 'E47'; Value = '  -3.37%  ' }
    @{ Cell = 'D48'; Value = '53.24' }
    @{ Cell = 'E48'; Value = '  -0.29%  ' }
    @{ Cell = 'D49'; Value = '72.36' }
    @{ Cell = 'E49'; Value = '  -1.10%  ' }
    @{ Cell = 'E50'; Value = '  +0.45%  ' }
    @{ Cell = 'D51'; Value = '90.99' }
    @{ Cell = 'E51'; Value = '  -0.90%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
